# feat: add 2022-Q3 data
#
# Insert a brand-new "2022-Q3" worksheet (fund-holdings detail, same shape as
# the other quarterly sheets) right before "2022-Q2", and add the
# corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the "2022-Q3" sheet immediately before "2022-Q2".
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (all text).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q3Sheet.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 2]
}

# Data rows: A=index(number), B..F=text, G=text (except last row: number), H=number.
# Columns: idx, code, name, scale, stockPosition, positionRatio, marketValue, positionRank
$q3Rows = @(
    @(0, "004895", "华商鑫安灵活配置混合",           "2.11", "92.54", "4.44", "0.0937", 2),
    @(1, "159851", "华宝中证金融科技主题ETF",         "1.94", "98.27", "4.00", "0.0776", 5),
    @(2, "460009", "华泰柏瑞量化先行混合A",           "4.22", "93.06", "1.04", "0.0439", 8),
    @(3, "159804", "国寿安保国证创业板中盘精选88ETF", "1.10", "98.91", "1.92", "0.0211", 7),
    @(4, "516100", "华夏中证金融科技主题ETF",         "0.51", "96.79", "3.92", "0.0200", 5),
    @(5, "516860", "博时中证金融科技主题ETF",         "0.34", "98.57", "3.98", "0.0135", 5),
    @(6, "010246", "华泰柏瑞量化先行混合C",           "0.25", "93.06", "1.04", "0.0026", 8),
    @(7, "006942", "华泰柏瑞量化明选混合A",           "0.30", "90.61", "0.86", "0.0026", 9),
    @(8, "006943", "华泰柏瑞量化明选混合C",           "0.00", "90.61", "0.86", "0",      9)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3Sheet.Cells.Item($r, 1).Value = $row[0]

    $cB = $q3Sheet.Cells.Item($r, 2); $cB.NumberFormat = "@"; $cB.Value = $row[1]
    $cC = $q3Sheet.Cells.Item($r, 3); $cC.NumberFormat = "@"; $cC.Value = $row[2]
    $cD = $q3Sheet.Cells.Item($r, 4); $cD.NumberFormat = "@"; $cD.Value = $row[3]
    $cE = $q3Sheet.Cells.Item($r, 5); $cE.NumberFormat = "@"; $cE.Value = $row[4]
    $cF = $q3Sheet.Cells.Item($r, 6); $cF.NumberFormat = "@"; $cF.Value = $row[5]

    if ($r -eq 10) {
        # Last data row stores the market-value column as a real number (0),
        # matching the source data - everyone else keeps it as text.
        $q3Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $cG = $q3Sheet.Cells.Item($r, 7); $cG.NumberFormat = "@"; $cG.Value = $row[6]
    }

    $q3Sheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a row for 2022-Q3 at the top of
#    the data (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 0.28

# Re-number the index column (A) for the rows that shifted down one slot.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5
